$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Title value (row 5: Property "Title")
$ws.Range("B5").Value = "Codes pour caractériser la population de l'étude"

# Date value (row 8: Property "Date")
$ws.Range("B8").Value = "2023-10-19T15:25:12+00:00"

# Description value (row 11: Property "Description") - previously duplicated the Title text,
# now gets its own proper description text
$ws.Range("B11").Value = "Codes pour caractériser la population ciblée par l'étude"
